# [ADR-2614] update alns requirements
# Re-applies the regenerated ALNS shift-schedule solution grid onto the
# "Solution" sheet: each staff member's day-by-day shift code (DO/M1/A1/
# M2/A2/M3) is updated in place to match the new solver output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "M1"
$ws.Range("H2").Value = "M3"
$ws.Range("I2").Value = "M1"
$ws.Range("J2").Value = "DO"
$ws.Range("K2").Value = "M1"
$ws.Range("O2").Value = "M3"
$ws.Range("Q2").Value = "M1"
$ws.Range("T2").Value = "M3"
$ws.Range("X2").Value = "M1"
$ws.Range("Y2").Value = "M1"
$ws.Range("Z2").Value = "M1"
$ws.Range("AA2").Value = "M3"
$ws.Range("AB2").Value = "A1"

# Row 3
$ws.Range("E3").Value = "M1"
$ws.Range("F3").Value = "M2"
$ws.Range("G3").Value = "A1"
$ws.Range("H3").Value = "M3"
$ws.Range("I3").Value = "DO"
$ws.Range("K3").Value = "M1"
$ws.Range("L3").Value = "M1"
$ws.Range("N3").Value = "A2"
$ws.Range("O3").Value = "M3"
$ws.Range("P3").Value = "M1"
$ws.Range("Q3").Value = "DO"
$ws.Range("R3").Value = "M1"
$ws.Range("T3").Value = "M2"
$ws.Range("V3").Value = "M3"
$ws.Range("W3").Value = "M1"
$ws.Range("X3").Value = "DO"
$ws.Range("Y3").Value = "M2"
$ws.Range("AC3").Value = "M3"

# Row 4
$ws.Range("E4").Value = "M1"
$ws.Range("F4").Value = "M3"
$ws.Range("W4").Value = "M1"
$ws.Range("X4").Value = "DO"

# Row 5
$ws.Range("B5").Value = "M1"
$ws.Range("C5").Value = "DO"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "A1"
$ws.Range("J5").Value = "M2"
$ws.Range("K5").Value = "DO"
$ws.Range("M5").Value = "M3"
$ws.Range("N5").Value = "M3"
$ws.Range("O5").Value = "A1"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "M3"
$ws.Range("R5").Value = "M2"
$ws.Range("U5").Value = "M1"
$ws.Range("V5").Value = "A2"
$ws.Range("Y5").Value = "DO"
$ws.Range("AB5").Value = "M3"
$ws.Range("AC5").Value = "A1"

# Row 6
$ws.Range("B6").Value = "M3"
$ws.Range("C6").Value = "A1"
$ws.Range("D6").Value = "DO"
$ws.Range("E6").Value = "A1"
$ws.Range("F6").Value = "A1"
$ws.Range("G6").Value = "M1"
$ws.Range("H6").Value = "A1"
$ws.Range("J6").Value = "A2"
$ws.Range("K6").Value = "A2"
$ws.Range("L6").Value = "A1"
$ws.Range("M6").Value = "A2"
$ws.Range("N6").Value = "M3"
$ws.Range("O6").Value = "A1"
$ws.Range("Q6").Value = "A1"
$ws.Range("R6").Value = "A1"
$ws.Range("S6").Value = "A2"
$ws.Range("T6").Value = "A2"
$ws.Range("U6").Value = "M3"
$ws.Range("V6").Value = "A1"
$ws.Range("W6").Value = "A1"
$ws.Range("X6").Value = "A1"
$ws.Range("Y6").Value = "A2"
$ws.Range("Z6").Value = "DO"
$ws.Range("AA6").Value = "A2"
$ws.Range("AB6").Value = "M3"
$ws.Range("AC6").Value = "A2"

# Row 7
$ws.Range("C7").Value = "M3"
$ws.Range("D7").Value = "A1"
$ws.Range("E7").Value = "A1"
$ws.Range("F7").Value = "A1"
$ws.Range("H7").Value = "A1"
$ws.Range("I7").Value = "DO"
$ws.Range("J7").Value = "M3"
$ws.Range("K7").Value = "A1"
$ws.Range("L7").Value = "A2"
$ws.Range("M7").Value = "A1"
$ws.Range("N7").Value = "A2"
$ws.Range("O7").Value = "A2"
$ws.Range("P7").Value = "M3"
$ws.Range("Q7").Value = "A1"
$ws.Range("R7").Value = "DO"
$ws.Range("S7").Value = "A1"
$ws.Range("T7").Value = "A1"
$ws.Range("U7").Value = "A1"
$ws.Range("V7").Value = "A1"
$ws.Range("W7").Value = "DO"
$ws.Range("X7").Value = "M3"
$ws.Range("Y7").Value = "A1"
$ws.Range("Z7").Value = "A1"
$ws.Range("AA7").Value = "A1"
$ws.Range("AC7").Value = "A1"

# Row 8
$ws.Range("B8").Value = "M3"
$ws.Range("C8").Value = "A2"
$ws.Range("D8").Value = "A1"
$ws.Range("E8").Value = "A1"
$ws.Range("F8").Value = "A2"
$ws.Range("G8").Value = "M2"
$ws.Range("H8").Value = "DO"
$ws.Range("I8").Value = "A1"
$ws.Range("J8").Value = "A1"
$ws.Range("K8").Value = "A2"
$ws.Range("L8").Value = "A2"
$ws.Range("M8").Value = "A2"
$ws.Range("N8").Value = "M3"
$ws.Range("O8").Value = "DO"
$ws.Range("P8").Value = "M3"
$ws.Range("Q8").Value = "A2"
$ws.Range("R8").Value = "A2"
$ws.Range("S8").Value = "A1"
$ws.Range("T8").Value = "A2"
$ws.Range("U8").Value = "M2"
$ws.Range("V8").Value = "DO"
$ws.Range("W8").Value = "DO"
$ws.Range("X8").Value = "A1"
$ws.Range("Y8").Value = "A2"
$ws.Range("Z8").Value = "A2"
$ws.Range("AA8").Value = "A2"
$ws.Range("AB8").Value = "M3"

# Row 9
$ws.Range("B9").Value = "M3"
$ws.Range("C9").Value = "DO"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = "M1"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A2"
$ws.Range("H9").Value = "A1"
$ws.Range("K9").Value = "M3"
$ws.Range("L9").Value = "DO"
$ws.Range("M9").Value = "M2"
$ws.Range("N9").Value = "A1"
$ws.Range("O9").Value = "A1"
$ws.Range("V9").Value = "A2"
$ws.Range("X9").Value = "M1"
$ws.Range("Y9").Value = "M2"
$ws.Range("Z9").Value = "M3"
$ws.Range("AC9").Value = "A1"

# Row 10
$ws.Range("C10").Value = "M1"
$ws.Range("D10").Value = "A1"
$ws.Range("E10").Value = "A1"
$ws.Range("H10").Value = "M2"
$ws.Range("I10").Value = "A1"
$ws.Range("J10").Value = "DO"
$ws.Range("K10").Value = "M2"
$ws.Range("L10").Value = "M3"
$ws.Range("M10").Value = "M1"
$ws.Range("N10").Value = "A1"
$ws.Range("O10").Value = "M1"
$ws.Range("P10").Value = "M3"
$ws.Range("Q10").Value = "M2"
$ws.Range("R10").Value = "A2"
$ws.Range("S10").Value = "DO"
$ws.Range("T10").Value = "A2"
$ws.Range("U10").Value = "M2"
$ws.Range("V10").Value = "M2"
$ws.Range("X10").Value = "M1"
$ws.Range("Z10").Value = "M3"
$ws.Range("AB10").Value = "A2"
$ws.Range("AC10").Value = "M2"
